$d = $word.ActiveDocument

# Locate the paragraph that currently reads:
#   "2. System exibe mensagem 'A quantidade informada deve ser maior ou
#    igual a 01 (um)!' (MSG002) "
# inside the "AF[7] - Quantidade menor ou igual a zero" alternative flow.
$targetText = "A quantidade informada deve ser maior ou igual a 01 (um)!"
$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.Contains($targetText)) {
        $targetParagraph = $candidate
        break
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not locate the AF[7] step 2 paragraph"
}

$r = $targetParagraph.Range

# Build the replacement content as raw WordprocessingML: the old single
# paragraph (item "2.") is turned into three sequential steps (2, 3, 4),
# each one following the same run layout already used by its sibling
# step-paragraphs in this document (an empty leading run, the stray
# bookmarkEnd marker carried over from the flow's opening bookmark, then
# the visible text run).
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve"></w:t></w:r><w:bookmarkEnd w:id="2"/><w:r><w:rPr/><w:t xml:space="preserve">2. System registra a quantidade informada </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve"></w:t></w:r><w:bookmarkEnd w:id="2"/><w:r><w:rPr/><w:t xml:space="preserve">3. Usu&#225;rio do Sistema clica no bot&#227;o 'Calcular Desconto!' </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve"></w:t></w:r><w:bookmarkEnd w:id="2"/><w:r><w:rPr/><w:t xml:space="preserve">4. System exibe mensagem 'A quantidade informada deve ser maior ou igual a 01 (um)!' (MSG002) </w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
